$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.468.75"
Set-TextValue $ws.Range("E2") "  +3.85%  "

Set-TextValue $ws.Range("D3") "1.816.53"
Set-TextValue $ws.Range("E3") "  +5.12%  "

Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.50%  "

Set-TextValue $ws.Range("D5") "343.71"
Set-TextValue $ws.Range("E5") "  +3.12%  "

Set-TextValue $ws.Range("D6") "0.9995"
Set-TextValue $ws.Range("E6") "  +0.50%  "

Set-TextValue $ws.Range("D7") "0.3831"
Set-TextValue $ws.Range("E7") "  +3.26%  "

Set-TextValue $ws.Range("D8") "0.3514"
Set-TextValue $ws.Range("E8") "  +4.13%  "

Set-TextValue $ws.Range("D9") "49.00"
Set-TextValue $ws.Range("E9") "  -0.17%  "

Set-TextValue $ws.Range("E10") "  +2.49%  "

Set-TextValue $ws.Range("D11") "0.07823"
Set-TextValue $ws.Range("E11") "  +4.25%  "

Set-TextValue $ws.Range("D12") "1.001"
Set-TextValue $ws.Range("E12") "  +0.74%  "

Set-TextValue $ws.Range("E13") "  +8.84%  "

Set-TextValue $ws.Range("D14") "6.597"
Set-TextValue $ws.Range("E14") "  +3.55%  "

Set-TextValue $ws.Range("D15") "1.817.92"
Set-TextValue $ws.Range("E15") "  +6.00%  "

Set-TextValue $ws.Range("D16") "7.235"
Set-TextValue $ws.Range("E16") "  +3.73%  "

Set-TextValue $ws.Range("D17") "0.00001121"
Set-TextValue $ws.Range("E17") "  +3.18%  "

Set-TextValue $ws.Range("D18") "0.06728"
Set-TextValue $ws.Range("E18") "  +0.72%  "

Set-TextValue $ws.Range("D19") "86.33"
Set-TextValue $ws.Range("E19") "  +4.20%  "

Set-TextValue $ws.Range("D20") "0.9995"
Set-TextValue $ws.Range("E20") "  +0.41%  "

Set-TextValue $ws.Range("D21") "17.66"
Set-TextValue $ws.Range("E21") "  +5.92%  "

Set-TextValue $ws.Range("D22") "6.578"
Set-TextValue $ws.Range("E22") "  +6.86%  "

Set-TextValue $ws.Range("D23") "13.17"
Set-TextValue $ws.Range("E23") "  +0.06%  "

Set-TextValue $ws.Range("D24") "27.482.11"
Set-TextValue $ws.Range("E24") "  +4.10%  "

Set-TextValue $ws.Range("D25") "2.460"
Set-TextValue $ws.Range("E25") "  +0.06%  "

Set-TextValue $ws.Range("D26") "2.678"
Set-TextValue $ws.Range("E26") "  +6.87%  "

Set-TextValue $ws.Range("D27") "22.19"
Set-TextValue $ws.Range("E27") "  +13.64%  "

Set-TextValue $ws.Range("D28") "1.469"
Set-TextValue $ws.Range("E28") "  +1.36%  "

Set-TextValue $ws.Range("D29") "154.01"
Set-TextValue $ws.Range("E29") "  +1.49%  "

Set-TextValue $ws.Range("D30") "2.023.62"
Set-TextValue $ws.Range("E30") "  +6.15%  "

Set-TextValue $ws.Range("D31") "136.38"
Set-TextValue $ws.Range("E31") "  +3.98%  "

Set-TextValue $ws.Range("D32") "6.350"
Set-TextValue $ws.Range("E32") "  +4.39%  "

Set-TextValue $ws.Range("D33") "4.064"
Set-TextValue $ws.Range("E33") "  -1.07%  "

Set-TextValue $ws.Range("D34") "13.89"
Set-TextValue $ws.Range("E34") "  +5.62%  "

Set-TextValue $ws.Range("D35") "0.08803"
Set-TextValue $ws.Range("E35") "  +2.77%  "

Set-TextValue $ws.Range("D36") "1.689"
Set-TextValue $ws.Range("E36") "  -1.50%  "

Set-TextValue $ws.Range("D37") "5.632"
Set-TextValue $ws.Range("E37") "  +3.62%  "

Set-TextValue $ws.Range("D38") "0.7002"
Set-TextValue $ws.Range("E38") "  +12.11%  "

Set-TextValue $ws.Range("D39") "0.2265"
Set-TextValue $ws.Range("E39") "  +4.89%  "

Set-TextValue $ws.Range("D40") "0.02408"
Set-TextValue $ws.Range("E40") "  +2.28%  "

Set-TextValue $ws.Range("D41") "0.06483"
Set-TextValue $ws.Range("E41") "  +2.52%  "

Set-TextValue $ws.Range("D42") "8.964"
Set-TextValue $ws.Range("E42") "  +3.46%  "

Set-TextValue $ws.Range("D43") "1.303"
Set-TextValue $ws.Range("E43") "  +5.26%  "

Set-TextValue $ws.Range("D44") "14.82"
Set-TextValue $ws.Range("E44") "  +3.10%  "

Set-TextValue $ws.Range("D45") "0.6575"
Set-TextValue $ws.Range("E45") "  +9.09%  "

Set-TextValue $ws.Range("D46") "0.9993"
Set-TextValue $ws.Range("E46") "  +0.39%  "

Set-TextValue $ws.Range("D47") "3.961"
Set-TextValue $ws.Range("E47") "  +1.70%  "

Set-TextValue $ws.Range("D48") "2.183"
Set-TextValue $ws.Range("E48") "  +6.40%  "

Set-TextValue $ws.Range("D49") "132.74"
Set-TextValue $ws.Range("E49") "  +2.64%  "

Set-TextValue $ws.Range("E50") "  -0.13%  "

Set-TextValue $ws.Range("D51") "80.55"
Set-TextValue $ws.Range("E51") "  +3.89%  "
